$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44235
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 14000
$ws.Range("P2").Value = 778

# Row 3
$ws.Range("D3").Value = 44235
$ws.Range("J3").Value = 70
$ws.Range("K3").Value = 12000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 12000
$ws.Range("P3").Value = 667

# Row 4
$ws.Range("D4").Value = 44235
$ws.Range("I4").Value = "Tercera"
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = 10000
$ws.Range("P4").Value = 556

# Row 5
$ws.Range("D5").Value = 44238
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 90
$ws.Range("K5").Value = 13000
$ws.Range("L5").Value = 13000
$ws.Range("M5").Value = 13000
$ws.Range("P5").Value = 722

# Row 6
$ws.Range("D6").Value = 44238
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 80
$ws.Range("K6").Value = 11000
$ws.Range("L6").Value = 11000
$ws.Range("M6").Value = 11000
$ws.Range("P6").Value = 611

# Row 7
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 75
$ws.Range("K7").Value = 13000
$ws.Range("L7").Value = 13000
$ws.Range("M7").Value = 13000
$ws.Range("P7").Value = 722

# Row 8
$ws.Range("D8").Value = 44991
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 56
$ws.Range("K8").Value = 9000
$ws.Range("L8").Value = 9000
$ws.Range("M8").Value = 9000
$ws.Range("P8").Value = 500

# Row 9
$ws.Range("D9").Value = 44424
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 75
$ws.Range("K9").Value = 18000
$ws.Range("L9").Value = 18000
$ws.Range("M9").Value = 18000
$ws.Range("N9").Value = "$/caja 15 kilos"
$ws.Range("P9").Value = 1200
$ws.Range("Q9").Value = 15

# Row 10
$ws.Range("D10").Value = 44424
$ws.Range("I10").Value = "Segunda"
$ws.Range("J10").Value = 50
$ws.Range("K10").Value = 12000
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 12000
$ws.Range("P10").Value = 800

# Row 11
$ws.Range("D11").Value = 44992
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 56
$ws.Range("K11").Value = 13000
$ws.Range("L11").Value = 13000
$ws.Range("M11").Value = 13000
$ws.Range("N11").Value = "$/bandeja 18 kilos"
$ws.Range("P11").Value = 722
$ws.Range("Q11").Value = 18

# Row 12
$ws.Range("D12").Value = 44756
$ws.Range("J12").Value = 65
$ws.Range("N12").Value = "$/caja 15 kilos"
$ws.Range("P12").Value = 933
$ws.Range("Q12").Value = 15

# Row 13
$ws.Range("D13").Value = 44756
$ws.Range("J13").Value = 68
$ws.Range("N13").Value = "$/caja 15 kilos"
$ws.Range("P13").Value = 800
$ws.Range("Q13").Value = 15

# Row 14
$ws.Range("D14").Value = 44242
$ws.Range("I14").Value = "Primera"
$ws.Range("K14").Value = 13000
$ws.Range("L14").Value = 13000
$ws.Range("M14").Value = 13000
$ws.Range("P14").Value = 722

# Row 15
$ws.Range("D15").Value = 44242
$ws.Range("I15").Value = "Segunda"
$ws.Range("J15").Value = 50
$ws.Range("K15").Value = 10000
$ws.Range("L15").Value = 10000
$ws.Range("M15").Value = 10000
$ws.Range("N15").Value = "$/bandeja 18 kilos"
$ws.Range("P15").Value = 556
$ws.Range("Q15").Value = 18

# Row 16
$ws.Range("D16").Value = 44536
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 87
$ws.Range("K16").Value = 22000
$ws.Range("L16").Value = 22000
$ws.Range("M16").Value = 22000
$ws.Range("N16").Value = "$/bandeja 18 kilos"
$ws.Range("P16").Value = 1222
$ws.Range("Q16").Value = 18

# Row 17
$ws.Range("D17").Value = 44536
$ws.Range("I17").Value = "Segunda"
$ws.Range("J17").Value = 80
$ws.Range("K17").Value = 20000
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = 20000
$ws.Range("P17").Value = 1111
